# Updates cryptos list data (prices & volume%) per upstream GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some Price (column D) values are numeric-looking strings (e.g. "57.00", "0.0785")
# that Excel would silently coerce to numbers (dropping trailing zeros / introducing
# floating point noise) if assigned directly. Force those specific cells to Text
# format first so the literal string is preserved exactly, matching the source feed.
$priceCellsAsText = @("D6", "D7", "D9", "D10", "D11", "D14", "D17", "D20", "D23", "D26", "D27", "D28", "D29", "D30", "D33", "D34", "D35", "D36", "D37", "D40", "D41", "D46")
foreach ($addr in $priceCellsAsText) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "36.501.60"
$ws.Range("E2").Value = "  +0.53%  "
$ws.Range("D3").Value = "2.007.06"
$ws.Range("E3").Value = "  +0.02%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "247.37"
$ws.Range("E5").Value = "  -1.66%  "
$ws.Range("D6").Value = "0.631"
$ws.Range("E6").Value = "  -1.12%  "
$ws.Range("D7").Value = "62.30"
$ws.Range("E7").Value = "  +1.47%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "0.384"
$ws.Range("E9").Value = "  +4.16%  "
$ws.Range("D10").Value = "57.00"
$ws.Range("E10").Value = "  -1.93%  "
$ws.Range("D11").Value = "0.0785"
$ws.Range("E11").Value = "  +6.16%  "
$ws.Range("E12").Value = "  +0.11%  "
$ws.Range("D13").Value = "0.879"
$ws.Range("E13").Value = "  -1.90%  "
$ws.Range("D14").Value = "22.69"
$ws.Range("E14").Value = "  +12.90%  "
$ws.Range("D15").Value = "14.07"
$ws.Range("E15").Value = "  -5.02%  "
$ws.Range("D16").Value = "2.301.23"
$ws.Range("E16").Value = "  -0.15%  "
$ws.Range("D17").Value = "5.52"
$ws.Range("E17").Value = "  +1.72%  "
$ws.Range("D18").Value = "2.009.07"
$ws.Range("E18").Value = "  -0.33%  "
$ws.Range("D19").Value = "36.455.54"
$ws.Range("E19").Value = "  +0.39%  "
$ws.Range("D20").Value = "71.91"
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("E21").Value = "  +1.51%  "
$ws.Range("E22").Value = "  +1.62%  "
$ws.Range("D23").Value = "238.48"
$ws.Range("E23").Value = "  +1.96%  "
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("D25").Value = "2.52"
$ws.Range("E25").Value = "  -6.94%  "
$ws.Range("D26").Value = "2.32"
$ws.Range("E26").Value = "  +0.32%  "
$ws.Range("D27").Value = "9.98"
$ws.Range("E27").Value = "  +4.66%  "
$ws.Range("D28").Value = "159.45"
$ws.Range("E28").Value = "  -2.36%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "20.12"
$ws.Range("E29").Value = "  +2.86%  "
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").Value = "0.134"
$ws.Range("E30").Value = "  +24.10%  "
$ws.Range("E31").Value = "  +1.13%  "
$ws.Range("E32").Value = "  -1.55%  "
$ws.Range("D33").Value = "1.17"
$ws.Range("E33").Value = "  -0.32%  "
$ws.Range("D34").Value = "0.0629"
$ws.Range("E34").Value = "  +3.94%  "
$ws.Range("D35").Value = "4.49"
$ws.Range("E35").Value = "  -1.27%  "
$ws.Range("D36").Value = "6.51"
$ws.Range("E36").Value = "  +10.47%  "
$ws.Range("D37").Value = "2.34"
$ws.Range("E37").Value = "  -4.00%  "
$ws.Range("E38").Value = "  +0.11%  "
$ws.Range("E39").Value = "  +1.06%  "
$ws.Range("D40").Value = "3.20"
$ws.Range("E40").Value = "  +12.51%  "
$ws.Range("B41").Value = "Cronos"
$ws.Range("C41").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D41").Value = "0.101"
$ws.Range("E41").Value = "  -2.98%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "1.26"
$ws.Range("E42").Value = "  +3.28%  "
$ws.Range("E43").Value = "  -1.04%  "
$ws.Range("E44").Value = "  -0.20%  "
$ws.Range("E45").Value = "  -0.64%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "95.99"
$ws.Range("E46").Value = "  +1.74%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "16.77"
$ws.Range("E47").Value = "  -0.84%  "
$ws.Range("E48").Value = "  -5.27%  "
$ws.Range("D49").Value = "1.365.45"
$ws.Range("E49").Value = "  -5.84%  "
$ws.Range("E50").Value = "  -1.33%  "
$ws.Range("D51").Value = "2.192.92"
$ws.Range("E51").Value = "  -0.13%  "
